$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Week Beginning 2023-11-13) - rewrite Commit Message and add Week Details / Carry Digit,
# keep Challenges (What/How) text but swap which column holds which message
$ws.Range("C4").Value = "Continued setup - added to a Git repository on GitHub so that collaboration between home and school account is easy"
$ws.Range("D4").Value = " - Set up a GitHub repo that is accessible by both home and school accounts"
$ws.Range("E4").Value = " - Gantt Chart`n - Timeline`n"
$ws.Range("F4").Value = "Onedrive not working for sharing between school and home accounts."
$ws.Range("G4").Value = "Using my school and home GitHub accounts instead."

# Row 5 (Week Beginning 2023-11-20) - fill in Timeline Stage and Commit Message
$ws.Range("B5").Value = "TBC"
$ws.Range("C5").Value = "EPQ Session was cancelled, week was busy. Did some listening. (R) "

# Update the view selection to reflect the newly active cell (E5) and reset horizontal scroll
$ws.Range("E5").Select()
